$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 519.25
$ws.Range("I6").Value = 420.53333
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 1261.59999
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -1149.59999
$ws.Range("N6").Value = -6224

$ws.Range("H32").Value = 1585.4286
$ws.Range("I32").Value = 718.4
$ws.Range("J32").Value = 2067.111
$ws.Range("K32").Value = 718.4
$ws.Range("L32").Value = 2067.111
$ws.Range("M32").Value = -392.4
$ws.Range("N32").Value = -2719.111

$ws.Range("H52").Value = 105000
$ws.Range("J52").Value = 105000
$ws.Range("L52").Value = 315000
$ws.Range("N52").Value = -315320

$ws.Range("H121").Value = 372.1579
$ws.Range("J121").Value = 372.1579
$ws.Range("L121").Value = 1116.4737
$ws.Range("N121").Value = -4610.4737

$ws.Range("H128").Value = 35499
$ws.Range("I128").Value = 20666.666
$ws.Range("J128").Value = 41855.715
$ws.Range("K128").Value = 20666.666
$ws.Range("L128").Value = 41855.715
$ws.Range("M128").Value = -15686.666
$ws.Range("N128").Value = -51815.715

$ws.Range("H129").Value = 1252.25
$ws.Range("I129").Value = 335.125
$ws.Range("J129").Value = 1348.7894
$ws.Range("K129").Value = 1005.375
$ws.Range("L129").Value = 4046.3682
$ws.Range("M129").Value = 3994.625
$ws.Range("N129").Value = -14046.3682

$ws.Range("H137").Value = 693104.0600000001
$ws.Range("I137").Value = 1766830
$ws.Range("K137").Value = 5300490
$ws.Range("M137").Value = -5297940

$ws.Range("H138").Value = 2506.7576
$ws.Range("I138").Value = 1384.9412
$ws.Range("J138").Value = 3698.6875
$ws.Range("K138").Value = 4154.8236
$ws.Range("L138").Value = 11096.0625
$ws.Range("M138").Value = 985.1764000000003
$ws.Range("N138").Value = -21376.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 802.75
$ws.Range("I2").Value = 755.5
$ws.Range("J2").Value = 850
$ws.Range("K2").Value = 755.5
$ws.Range("L2").Value = 850
$ws.Range("M2").Value = -642.5
$ws.Range("N2").Value = -1076

$ws.Range("H110").Value = 3346.6
$ws.Range("I110").Value = 3627.75
$ws.Range("K110").Value = 3627.75
$ws.Range("M110").Value = -1582.75

$ws.Range("H116").Value = 802.75
$ws.Range("I116").Value = 755.5
$ws.Range("J116").Value = 850
$ws.Range("K116").Value = 755.5
$ws.Range("L116").Value = 850
$ws.Range("M116").Value = 1538.5
$ws.Range("N116").Value = -5438

$ws.Range("H122").Value = 3451.625
$ws.Range("I122").Value = 1304
$ws.Range("K122").Value = 3912
$ws.Range("M122").Value = -1462

$ws.Range("H132").Value = 2291.7334
$ws.Range("I132").Value = 1765.6757
$ws.Range("J132").Value = 4724.75
$ws.Range("K132").Value = 5297.0271
$ws.Range("L132").Value = 14174.25
$ws.Range("M132").Value = -2767.0271
$ws.Range("N132").Value = -19234.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 802.75
$ws.Range("I3").Value = 755.5
$ws.Range("J3").Value = 850
$ws.Range("K3").Value = 755.5
$ws.Range("L3").Value = 850
$ws.Range("M3").Value = -641.5
$ws.Range("N3").Value = -1078

$ws.Range("H118").Value = 29890
$ws.Range("J118").Value = 29890
$ws.Range("L118").Value = 29890
$ws.Range("N118").Value = -33204

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 252702.72
$ws.Range("I31").Value = 451276.03
$ws.Range("J31").Value = 4486.0835
$ws.Range("K31").Value = 451276.03
$ws.Range("L31").Value = 4486.0835
$ws.Range("M31").Value = -450981.03
$ws.Range("N31").Value = -5076.0835

$ws.Range("H34").Value = 252702.72
$ws.Range("I34").Value = 451276.03
$ws.Range("J34").Value = 4486.0835
$ws.Range("K34").Value = 451276.03
$ws.Range("L34").Value = 4486.0835
$ws.Range("M34").Value = -451074.03
$ws.Range("N34").Value = -4890.0835

$ws.Range("H109").Value = 34886.332
$ws.Range("J109").Value = 34886.332
$ws.Range("L109").Value = 34886.332
$ws.Range("N109").Value = -36966.332

$ws.Range("H122").Value = 6952.3335
$ws.Range("J122").Value = 8957
$ws.Range("L122").Value = 26871
$ws.Range("N122").Value = -31771

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1232.409
$ws.Range("I68").Value = 1021.53125
$ws.Range("J68").Value = 1430.8823
$ws.Range("K68").Value = 3064.59375
$ws.Range("L68").Value = 4292.6469
$ws.Range("M68").Value = -2253.59375
$ws.Range("N68").Value = -5914.6469

$ws.Range("H71").Value = 1232.409
$ws.Range("I71").Value = 1021.53125
$ws.Range("J71").Value = 1430.8823
$ws.Range("K71").Value = 9193.78125
$ws.Range("L71").Value = 12877.9407
$ws.Range("M71").Value = -5137.78125
$ws.Range("N71").Value = -20989.9407

$ws.Range("H107").Value = 27681.842
$ws.Range("J107").Value = 94691.37
$ws.Range("L107").Value = 284074.11
$ws.Range("N107").Value = -287914.11

$ws.Range("H112").Value = 7370
$ws.Range("I112").Value = 425
$ws.Range("J112").Value = 12000
$ws.Range("K112").Value = 1275
$ws.Range("L112").Value = 36000
$ws.Range("M112").Value = -167
$ws.Range("N112").Value = -38216

$ws.Range("H113").Value = 1812155.8
$ws.Range("I113").Value = 570.6042
$ws.Range("J113").Value = 5952922
$ws.Range("K113").Value = 1711.8126
$ws.Range("L113").Value = 17858766
$ws.Range("M113").Value = 458.1874
$ws.Range("N113").Value = -17863106

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6533.6665
$ws.Range("I40").Value = 5475.375
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 5475.375
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -5339.375
$ws.Range("N40").Value = -15272

$ws.Range("H122").Value = 7118.1816
$ws.Range("I122").Value = 5500.25
$ws.Range("J122").Value = 8042.7144
$ws.Range("K122").Value = 16500.75
$ws.Range("L122").Value = 24128.1432
$ws.Range("M122").Value = -14050.75
$ws.Range("N122").Value = -29028.1432

$ws.Range("H136").Value = 2610.0217
$ws.Range("I136").Value = 1051.7587
$ws.Range("J136").Value = 5268.2354
$ws.Range("K136").Value = 3155.2761
$ws.Range("L136").Value = 15804.7062
$ws.Range("M136").Value = -605.2761
$ws.Range("N136").Value = -20904.7062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9799.799999999999
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 9799.799999999999
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 29399.4
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -34299.39999999999

$ws.Range("H126").Value = 1185057
$ws.Range("I126").Value = 1702.6666
$ws.Range("J126").Value = 3551765.8
$ws.Range("K126").Value = 5107.9998
$ws.Range("L126").Value = 10655297.4
$ws.Range("M126").Value = -2637.9998
$ws.Range("N126").Value = -10660237.4

